$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-02 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-03 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("30-26=", $true, $false, $false, $false, $false, $true, 1, $false, "75-30=", 2) | Out-Null
$d.Content.Find.Execute("17+67=", $true, $false, $false, $false, $false, $true, 1, $false, "76-13=", 2) | Out-Null
$d.Content.Find.Execute("14+77=", $true, $false, $false, $false, $false, $true, 1, $false, "78-11=", 2) | Out-Null
$d.Content.Find.Execute("93-79=", $true, $false, $false, $false, $false, $true, 1, $false, "75-70=", 2) | Out-Null
$d.Content.Find.Execute("33-9=", $true, $false, $false, $false, $false, $true, 1, $false, "37+28=", 2) | Out-Null
$d.Content.Find.Execute("95-24=", $true, $false, $false, $false, $false, $true, 1, $false, "31-28=", 2) | Out-Null
$d.Content.Find.Execute("45-6=", $true, $false, $false, $false, $false, $true, 1, $false, "56+38=", 2) | Out-Null
$d.Content.Find.Execute("58+19=", $true, $false, $false, $false, $false, $true, 1, $false, "30+36=", 2) | Out-Null
$d.Content.Find.Execute("28+54=", $true, $false, $false, $false, $false, $true, 1, $false, "37+49=", 2) | Out-Null
$d.Content.Find.Execute("58-19=", $true, $false, $false, $false, $false, $true, 1, $false, "82-72=", 2) | Out-Null
$d.Content.Find.Execute("20+68=", $true, $false, $false, $false, $false, $true, 1, $false, "45-27=", 2) | Out-Null
$d.Content.Find.Execute("47+47=", $true, $false, $false, $false, $false, $true, 1, $false, "73-56=", 2) | Out-Null
$d.Content.Find.Execute("35-1=", $true, $false, $false, $false, $false, $true, 1, $false, "12+57=", 2) | Out-Null
$d.Content.Find.Execute("87-42=", $true, $false, $false, $false, $false, $true, 1, $false, "18+30=", 2) | Out-Null
$d.Content.Find.Execute("94+5=", $true, $false, $false, $false, $false, $true, 1, $false, "61-5=", 2) | Out-Null
$d.Content.Find.Execute("61-7=", $true, $false, $false, $false, $false, $true, 1, $false, "43-2=", 2) | Out-Null
$d.Content.Find.Execute("33+38=", $true, $false, $false, $false, $false, $true, 1, $false, "89-72=", 2) | Out-Null
$d.Content.Find.Execute("55+25=", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=", 2) | Out-Null
$d.Content.Find.Execute("58-52=", $true, $false, $false, $false, $false, $true, 1, $false, "0+83=", 2) | Out-Null
$d.Content.Find.Execute("28+32=", $true, $false, $false, $false, $false, $true, 1, $false, "8+70=", 2) | Out-Null
$d.Content.Find.Execute("69+20=", $true, $false, $false, $false, $false, $true, 1, $false, "11-2=", 2) | Out-Null
$d.Content.Find.Execute("18+3=", $true, $false, $false, $false, $false, $true, 1, $false, "52-31=", 2) | Out-Null
$d.Content.Find.Execute("70+19=", $true, $false, $false, $false, $false, $true, 1, $false, "18+66=", 2) | Out-Null
$d.Content.Find.Execute("65-38=", $true, $false, $false, $false, $false, $true, 1, $false, "15+5=", 2) | Out-Null
$d.Content.Find.Execute("92-38=", $true, $false, $false, $false, $false, $true, 1, $false, "73-3=", 2) | Out-Null
$d.Content.Find.Execute("41+39=", $true, $false, $false, $false, $false, $true, 1, $false, "55-33=", 2) | Out-Null
$d.Content.Find.Execute("31+38=", $true, $false, $false, $false, $false, $true, 1, $false, "88-51=", 2) | Out-Null
$d.Content.Find.Execute("85-74=", $true, $false, $false, $false, $false, $true, 1, $false, "91-17=", 2) | Out-Null
$d.Content.Find.Execute("15-6=", $true, $false, $false, $false, $false, $true, 1, $false, "13+10=", 2) | Out-Null
$d.Content.Find.Execute("99-70=", $true, $false, $false, $false, $false, $true, 1, $false, "87+2=", 2) | Out-Null
$d.Content.Find.Execute("77-9=", $true, $false, $false, $false, $false, $true, 1, $false, "49+5=", 2) | Out-Null
$d.Content.Find.Execute("49+22=", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=", 2) | Out-Null
$d.Content.Find.Execute("29-27=", $true, $false, $false, $false, $false, $true, 1, $false, "17+40=", 2) | Out-Null
$d.Content.Find.Execute("37+20=", $true, $false, $false, $false, $false, $true, 1, $false, "50-12=", 2) | Out-Null
$d.Content.Find.Execute("69-29=", $true, $false, $false, $false, $false, $true, 1, $false, "85-72=", 2) | Out-Null
$d.Content.Find.Execute("83-60=", $true, $false, $false, $false, $false, $true, 1, $false, "10+0=", 2) | Out-Null
$d.Content.Find.Execute("34+57=", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=", 2) | Out-Null
$d.Content.Find.Execute("7+12=", $true, $false, $false, $false, $false, $true, 1, $false, "45+52=", 2) | Out-Null
$d.Content.Find.Execute("48-16=", $true, $false, $false, $false, $false, $true, 1, $false, "65-8=", 2) | Out-Null
$d.Content.Find.Execute("30-8=", $true, $false, $false, $false, $false, $true, 1, $false, "90-50=", 2) | Out-Null
$d.Content.Find.Execute("27+26=", $true, $false, $false, $false, $false, $true, 1, $false, "21+70=", 2) | Out-Null
$d.Content.Find.Execute("27+10=", $true, $false, $false, $false, $false, $true, 1, $false, "27+62=", 2) | Out-Null
$d.Content.Find.Execute("6+83=", $true, $false, $false, $false, $false, $true, 1, $false, "49+27=", 2) | Out-Null
$d.Content.Find.Execute("63+28=", $true, $false, $false, $false, $false, $true, 1, $false, "57-51=", 2) | Out-Null
$d.Content.Find.Execute("47+43=", $true, $false, $false, $false, $false, $true, 1, $false, "57+12=", 2) | Out-Null
$d.Content.Find.Execute("34-6=", $true, $false, $false, $false, $false, $true, 1, $false, "25+31=", 2) | Out-Null
$d.Content.Find.Execute("74-17=", $true, $false, $false, $false, $false, $true, 1, $false, "45+17=", 2) | Out-Null
$d.Content.Find.Execute("38+25=", $true, $false, $false, $false, $false, $true, 1, $false, "54+45=", 2) | Out-Null
$d.Content.Find.Execute("80-51=", $true, $false, $false, $false, $false, $true, 1, $false, "22+7=", 2) | Out-Null
$d.Content.Find.Execute("89-61=", $true, $false, $false, $false, $false, $true, 1, $false, "82-53=", 2) | Out-Null
$d.Content.Find.Execute("94-25=", $true, $false, $false, $false, $false, $true, 1, $false, "55+29=", 2) | Out-Null
$d.Content.Find.Execute("0+43=", $true, $false, $false, $false, $false, $true, 1, $false, "85-39=", 2) | Out-Null
$d.Content.Find.Execute("56+29=", $true, $false, $false, $false, $false, $true, 1, $false, "87-80=", 2) | Out-Null
$d.Content.Find.Execute("57-42=", $true, $false, $false, $false, $false, $true, 1, $false, "82-12=", 2) | Out-Null
$d.Content.Find.Execute("40+34=", $true, $false, $false, $false, $false, $true, 1, $false, "4+1=", 2) | Out-Null
$d.Content.Find.Execute("79+7=", $true, $false, $false, $false, $false, $true, 1, $false, "22-14=", 2) | Out-Null
$d.Content.Find.Execute("6-2=", $true, $false, $false, $false, $false, $true, 1, $false, "77-12=", 2) | Out-Null
$d.Content.Find.Execute("51+44=", $true, $false, $false, $false, $false, $true, 1, $false, "55+7=", 2) | Out-Null
$d.Content.Find.Execute("36+25=", $true, $false, $false, $false, $false, $true, 1, $false, "79-35=", 2) | Out-Null
$d.Content.Find.Execute("4+27=", $true, $false, $false, $false, $false, $true, 1, $false, "37+53=", 2) | Out-Null
$d.Content.Find.Execute("48+18=", $true, $false, $false, $false, $false, $true, 1, $false, "38+57=", 2) | Out-Null
$d.Content.Find.Execute("51-35=", $true, $false, $false, $false, $false, $true, 1, $false, "47+1=", 2) | Out-Null
$d.Content.Find.Execute("64-38=", $true, $false, $false, $false, $false, $true, 1, $false, "8+39=", 2) | Out-Null
$d.Content.Find.Execute("91+1=", $true, $false, $false, $false, $false, $true, 1, $false, "80-17=", 2) | Out-Null
$d.Content.Find.Execute("94-28=", $true, $false, $false, $false, $false, $true, 1, $false, "39-24=", 2) | Out-Null
$d.Content.Find.Execute("68-0=", $true, $false, $false, $false, $false, $true, 1, $false, "46-45=", 2) | Out-Null
$d.Content.Find.Execute("84-9=", $true, $false, $false, $false, $false, $true, 1, $false, "0+82=", 2) | Out-Null
$d.Content.Find.Execute("91-44=", $true, $false, $false, $false, $false, $true, 1, $false, "68-17=", 2) | Out-Null
$d.Content.Find.Execute("61-57=", $true, $false, $false, $false, $false, $true, 1, $false, "58-30=", 2) | Out-Null
$d.Content.Find.Execute("60-2=", $true, $false, $false, $false, $false, $true, 1, $false, "84-51=", 2) | Out-Null
$d.Content.Find.Execute("56+23=", $true, $false, $false, $false, $false, $true, 1, $false, "10+31=", 2) | Out-Null
$d.Content.Find.Execute("82+7=", $true, $false, $false, $false, $false, $true, 1, $false, "22+19=", 2) | Out-Null
$d.Content.Find.Execute("28+44=", $true, $false, $false, $false, $false, $true, 1, $false, "37-16=", 2) | Out-Null
$d.Content.Find.Execute("18+60=", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=", 2) | Out-Null
$d.Content.Find.Execute("56+30=", $true, $false, $false, $false, $false, $true, 1, $false, "39-36=", 2) | Out-Null
$d.Content.Find.Execute("77-69=", $true, $false, $false, $false, $false, $true, 1, $false, "43-2=", 2) | Out-Null
$d.Content.Find.Execute("56-48=", $true, $false, $false, $false, $false, $true, 1, $false, "7+18=", 2) | Out-Null
$d.Content.Find.Execute("90-79=", $true, $false, $false, $false, $false, $true, 1, $false, "38+45=", 2) | Out-Null
$d.Content.Find.Execute("75-54=", $true, $false, $false, $false, $false, $true, 1, $false, "20-3=", 2) | Out-Null
$d.Content.Find.Execute("14+64=", $true, $false, $false, $false, $false, $true, 1, $false, "54-33=", 2) | Out-Null
$d.Content.Find.Execute("79-66=", $true, $false, $false, $false, $false, $true, 1, $false, "66-20=", 2) | Out-Null
$d.Content.Find.Execute("86-7=", $true, $false, $false, $false, $false, $true, 1, $false, "69+5=", 2) | Out-Null
$d.Content.Find.Execute("1+10=", $true, $false, $false, $false, $false, $true, 1, $false, "61+34=", 2) | Out-Null
$d.Content.Find.Execute("69+15=", $true, $false, $false, $false, $false, $true, 1, $false, "32-21=", 2) | Out-Null
$d.Content.Find.Execute("50+22=", $true, $false, $false, $false, $false, $true, 1, $false, "1+25=", 2) | Out-Null
$d.Content.Find.Execute("27+9=", $true, $false, $false, $false, $false, $true, 1, $false, "95-77=", 2) | Out-Null
$d.Content.Find.Execute("34+29=", $true, $false, $false, $false, $false, $true, 1, $false, "42+33=", 2) | Out-Null
$d.Content.Find.Execute("26+4=", $true, $false, $false, $false, $false, $true, 1, $false, "28+61=", 2) | Out-Null
$d.Content.Find.Execute("24-9=", $true, $false, $false, $false, $false, $true, 1, $false, "81+12=", 2) | Out-Null
$d.Content.Find.Execute("77+11=", $true, $false, $false, $false, $false, $true, 1, $false, "42-40=", 2) | Out-Null
$d.Content.Find.Execute("5+77=", $true, $false, $false, $false, $false, $true, 1, $false, "5+65=", 2) | Out-Null
$d.Content.Find.Execute("54+44=", $true, $false, $false, $false, $false, $true, 1, $false, "81-25=", 2) | Out-Null
$d.Content.Find.Execute("62+27=", $true, $false, $false, $false, $false, $true, 1, $false, "95-73=", 2) | Out-Null
$d.Content.Find.Execute("50-18=", $true, $false, $false, $false, $false, $true, 1, $false, "54-45=", 2) | Out-Null
$d.Content.Find.Execute("60+18=", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=", 2) | Out-Null
$d.Content.Find.Execute("11+72=", $true, $false, $false, $false, $false, $true, 1, $false, "18+15=", 2) | Out-Null
$d.Content.Find.Execute("3+88=", $true, $false, $false, $false, $false, $true, 1, $false, "82-53=", 2) | Out-Null
$d.Content.Find.Execute("19+54=", $true, $false, $false, $false, $false, $true, 1, $false, "12-12=", 2) | Out-Null
$d.Content.Find.Execute("64-37=", $true, $false, $false, $false, $false, $true, 1, $false, "68+0=", 2) | Out-Null
$d.Content.Find.Execute("96-60=", $true, $false, $false, $false, $false, $true, 1, $false, "23+55=", 2) | Out-Null
